# Apply "contingencies with rene fine" edit:
#  - Renumbers the line/extr contingency rows so that line1..line8 and
#    extr1..extr8 are listed in one continuous sequence (16 data rows
#    instead of 14): rows 8-9 become line7/line8, rows 10-15 shift to
#    extr1..extr6, and two brand-new rows (16-17) hold extr7/extr8.
#  - Updates the C/D/E (from_bus/to_bus/in_service) values for every
#    row from row 8 down to the newly appended row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# id, name, from_bus(C), to_bus(D), in_service(E)
$rows = @(
    @{ R = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ R = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ R = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ R = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ R = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $true  },
    @{ R = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ R = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ R = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ R = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ R = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

# The two new rows (16-17) need the same cell formatting (bold, centred,
# boxed) that column A already carries on every data row; copy it from
# row 7 (an existing, untouched row) before writing the new values.
$ws.Cells.Item(7, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(7, 1).Copy($ws.Cells.Item(17, 1))

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}
